# Update scripts with new TPM values.
# The sending cluster "ECs" rows (originally rows 2-4) are removed, and the
# "MuSCs" sending-cluster rows (originally rows 5-7) move up to rows 2-4 with
# refreshed TPM-derived specificity values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three "ECs" sending-cluster rows (old rows 2,3,4).
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# Rewrite the now-shifted rows (2,3,4; were 5,6,7) with the refreshed TPM values.
$data = @(
    @("MuSCs","Cdh1","Egfr","ECs",3,1,0.9477166666666667,2.84315,1,1,3,1,0.4102596666666667,1.230779,0.003499619873322347,0.003499619873322347,0.3888099237611111,3.49928931385,0.003499619873322347,0.003499619873322347),
    @("MuSCs","Cdh1","Egfr","FAPs",3,1,0.9477166666666667,2.84315,1,1,3,1,101.898173,305.694519,0.8692174743460166,0.8692174743460165,96.57059685498334,869.1353716948501,0.8692174743460166,0.8692174743460165),
    @("MuSCs","Cdh1","Egfr","MuSCs",3,1,0.9477166666666667,2.84315,1,1,3,1,14.921347,44.764041,0.1272829057806611,0.1272829057806611,14.14120924101667,127.27088316915,0.1272829057806611,0.1272829057806611)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowValues[$c]
    }
}

Write-Host "done"
